$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.543.35'
$ws.Range("E2").Value = '  +3.07%  '
$ws.Range("D3").Value = '2.124.56'
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '346.99'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5247'
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4474'
$ws.Range("E8").Value = '  +1.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.43'
$ws.Range("E9").Value = '  +4.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09413'
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.183'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.31'
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.736'
$ws.Range("E13").Value = '  +6.42%  '
$ws.Range("D14").Value = '2.143.36'
$ws.Range("E14").Value = '  +3.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.985'
$ws.Range("E15").Value = '  +3.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '102.05'
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001169'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.008'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.53'
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06735'
$ws.Range("E20").Value = '  +1.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.367'
$ws.Range("E21").Value = '  +2.61%  '
$ws.Range("D23").Value = '30.597.68'
$ws.Range("E23").Value = '  +3.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.77'
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.333'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").Value = '2.381.04'
$ws.Range("E26").Value = '  +2.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.22'
$ws.Range("E27").Value = '  +1.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.556'
$ws.Range("E28").Value = '  +1.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.87'
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.32'
$ws.Range("E30").Value = '  +1.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.165'
$ws.Range("E31").Value = '  +1.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.790'
$ws.Range("E32").Value = '  +10.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.1063'
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.894'
$ws.Range("E34").Value = '  +12.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.309'
$ws.Range("E35").Value = '  +1.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.963'
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.56'
$ws.Range("E37").Value = '  +2.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02661'
$ws.Range("E38").Value = '  +3.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06877'
$ws.Range("E39").Value = '  +2.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.7156'
$ws.Range("E40").Value = '  +4.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.70'
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2256'
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.338'
$ws.Range("E43").Value = '  +4.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6939'
$ws.Range("E44").Value = '  +4.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.72'
$ws.Range("E45").Value = '  +3.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.395'
$ws.Range("E46").Value = '  +3.39%  '
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.387'
$ws.Range("E48").Value = '  +18.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.664'
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.233'
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000344'
$ws.Range("E51").Value = '  +3.57%  '
